# Applies the "sim_3_50" hyperparameter-table formatting fixes:
#  - corrects the mislabeled "RMSE" headers on the Double_Layer sheet to "LogLoss"
#  - renames the "Activation funcs" label to "Activation functions"
#  - applies a 2-decimal ("0.00") number format to the numeric result grids
#  - leaves Double_Layer as the active/selected sheet with updated selections

$wb = $excel.ActiveWorkbook

$wsSingle = $wb.Worksheets.Item("Single_Layer")
$wsDouble = $wb.Worksheets.Item("Double_Layer")

# --- Text fixes -----------------------------------------------------------

# "Activation funcs" -> "Activation functions" row label on both sheets
# (Double_Layer has the label twice, once per stacked table)
$wsSingle.Range("A4").Value = "Activation functions"
$wsDouble.Range("A4").Value = "Activation functions"
$wsDouble.Range("A19").Value = "Activation functions"

# Double_Layer's per-column headers were mislabeled "RMSE" - they should read
# "LogLoss", matching Single_Layer and the rest of the workbook.
foreach ($col in @("D", "F", "H", "J", "L", "N", "P", "R")) {
    $wsDouble.Range($col + "6").Value = "LogLoss"
    $wsDouble.Range($col + "21").Value = "LogLoss"
}

# --- Number formatting ------------------------------------------------------

# Apply a 2-decimal-place number format to the numeric result tables.
$wsSingle.Range("C7:R11").NumberFormat = "0.00"

$wsDouble.Range("C7:R11").NumberFormat = "0.00"
$wsDouble.Range("C22:R26").NumberFormat = "0.00"

# --- Sheet views / selection -------------------------------------------------

# Single_Layer: no longer the selected tab, selection moves to H18
$wsSingle.Range("H18").Select()

# Double_Layer: becomes the active/selected tab, scrolled/selected at M15
$wsDouble.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 7
$aw.ScrollColumn = 1
$wsDouble.Range("M15").Select()
